$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 14 - LIDAR Power Supply: set Design = "done (needs work)", Location = "in"
$ws.Range("C14").Value = "done (needs work)"
$ws.Range("D14").Value = "in"

# Row 15 - cRIO Cap: set Design = "done"; Location changes from "known" to "in"
$ws.Range("C15").Value = "done"
$ws.Range("D15").Value = "in"

# Row 17 - Sabertooth: set Design = "done (needs work)"; Location changes from "known" to "in"
$ws.Range("C17").Value = "done (needs work)"
$ws.Range("D17").Value = "in"

# Row 20 - Router: set Design = "done", Location = "in"
$ws.Range("C20").Value = "done"
$ws.Range("D20").Value = "in"

# Row 21 - Relay: set Design = "done", Location = "in"
$ws.Range("C21").Value = "done"
$ws.Range("D21").Value = "in"

# Row 32 (new) - 13.8 v supply power supply entry
$ws.Range("A32").Value = "13.8 v supply"
$ws.Range("C32").Value = "done"
$ws.Range("D32").Value = "in"

# Reflect the cursor position left in the sheet after editing
$ws.Range("F10").Select() | Out-Null
